$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63, shifting the existing rows 63:73 down to 64:74
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new record's data
$ws.Cells.Item(63, 1).Value = 10
$ws.Cells.Item(63, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(63, 3).Value = "La Araucanía"
$ws.Cells.Item(63, 4).Value = 44694
$ws.Cells.Item(63, 5).Value = 9
$ws.Cells.Item(63, 6).Value = 100112035
$ws.Cells.Item(63, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 25
$ws.Cells.Item(63, 11).Value = 30000
$ws.Cells.Item(63, 12).Value = 30000
$ws.Cells.Item(63, 13).Value = 30000
$ws.Cells.Item(63, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(63, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(63, 16).Value = 3000
$ws.Cells.Item(63, 17).Value = 10
$ws.Cells.Item(63, 18).Value = "Hortaliza"
